$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1442
$ws.Range("I19").Value = 1582.5714
$ws.Range("J19").Value = 950
$ws.Range("K19").Value = 1582.5714
$ws.Range("L19").Value = 950
$ws.Range("M19").Value = -1407.5714
$ws.Range("N19").Value = -1300
$ws.Range("H33").Value = 7061.1763
$ws.Range("I33").Value = 7964.846
$ws.Range("K33").Value = 7964.846
$ws.Range("M33").Value = -7735.846
$ws.Range("H62").Value = 12225
$ws.Range("I62").Value = 17097.666
$ws.Range("J62").Value = 9976.076999999999
$ws.Range("K62").Value = 17097.666
$ws.Range("L62").Value = 9976.076999999999
$ws.Range("M62").Value = -16473.666
$ws.Range("N62").Value = -11224.077
$ws.Range("H65").Value = 12225
$ws.Range("I65").Value = 17097.666
$ws.Range("J65").Value = 9976.076999999999
$ws.Range("K65").Value = 85488.33
$ws.Range("L65").Value = 49880.38499999999
$ws.Range("M65").Value = -82368.33
$ws.Range("N65").Value = -56120.38499999999
$ws.Range("H74").Value = 3759.4
$ws.Range("I74").Value = 3519.8
$ws.Range("K74").Value = 3519.8
$ws.Range("M74").Value = -2583.8
$ws.Range("H77").Value = 3759.4
$ws.Range("I77").Value = 3519.8
$ws.Range("K77").Value = 17599
$ws.Range("M77").Value = -12919
$ws.Range("H92").Value = 46094
$ws.Range("I92").Value = 83818.75
$ws.Range("K92").Value = 83818.75
$ws.Range("M92").Value = -82570.75
$ws.Range("H97").Value = 2152.9167
$ws.Range("J97").Value = 2642.7778
$ws.Range("L97").Value = 7928.3334
$ws.Range("N97").Value = -8920.3334
$ws.Range("H98").Value = 1728.5
$ws.Range("I98").Value = 1940
$ws.Range("K98").Value = 1940
$ws.Range("M98").Value = -442
$ws.Range("H99").Value = 2494.5
$ws.Range("J99").Value = 5906.75
$ws.Range("L99").Value = 17720.25
$ws.Range("N99").Value = -20716.25
$ws.Range("H100").Value = 2705
$ws.Range("I100").Value = 2016.6666
$ws.Range("K100").Value = 2016.6666
$ws.Range("M100").Value = -1475.6666
$ws.Range("H106").Value = 10852.5
$ws.Range("I106").Value = 8812.5
$ws.Range("K106").Value = 8812.5
$ws.Range("M106").Value = -8181.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 5687.8
$ws.Range("J113").Value = 6064.143
$ws.Range("L113").Value = 6064.143
$ws.Range("N113").Value = -12572.143
$ws.Range("H115").Value = 300
$ws.Range("I115").Value = 300
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 900
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 667
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 5522.923
$ws.Range("I116").Value = 4930.1
$ws.Range("K116").Value = 4930.1
$ws.Range("M116").Value = -1488.1
$ws.Range("H118").Value = 598.1429000000001
$ws.Range("I118").Value = 360
$ws.Range("J118").Value = 1193.5
$ws.Range("K118").Value = 1080
$ws.Range("L118").Value = 3580.5
$ws.Range("M118").Value = 577
$ws.Range("N118").Value = -6894.5
$ws.Range("H122").Value = 1728.5
$ws.Range("I122").Value = 1940
$ws.Range("K122").Value = 5820
$ws.Range("M122").Value = -3370
$ws.Range("H131").Value = 2230.5293
$ws.Range("I131").Value = 1119.9375
$ws.Range("K131").Value = 3359.8125
$ws.Range("M131").Value = 1680.1875
$ws.Range("H137").Value = 2160.4
$ws.Range("I137").Value = 1134.3334
$ws.Range("J137").Value = 3699.5
$ws.Range("K137").Value = 3403.0002
$ws.Range("L137").Value = 11098.5
$ws.Range("M137").Value = -853.0001999999999
$ws.Range("N137").Value = -16198.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 18000
$ws.Range("I22").Value = 18000
$ws.Range("K22").Value = 18000
$ws.Range("M22").Value = -17701
$ws.Range("H43").Value = 59999.5
$ws.Range("J43").Value = 19999
$ws.Range("L43").Value = 19999
$ws.Range("N43").Value = -20625
$ws.Range("H45").Value = 3613.95
$ws.Range("I45").Value = 2500.8
$ws.Range("K45").Value = 2500.8
$ws.Range("M45").Value = -2123.8
$ws.Range("H61").Value = 1955
$ws.Range("I61").Value = 1786.3334
$ws.Range("K61").Value = 1786.3334
$ws.Range("M61").Value = -1574.3334
$ws.Range("H74").Value = 27020.309
$ws.Range("I74").Value = 29338.236
$ws.Range("K74").Value = 29338.236
$ws.Range("M74").Value = -28464.236
$ws.Range("H77").Value = 27020.309
$ws.Range("I77").Value = 29338.236
$ws.Range("K77").Value = 146691.18
$ws.Range("M77").Value = -142323.18
$ws.Range("H88").Value = 1152.4
$ws.Range("J88").Value = 1514.5454
$ws.Range("L88").Value = 1514.5454
$ws.Range("N88").Value = -2326.5454
$ws.Range("H91").Value = 1152.4
$ws.Range("J91").Value = 1514.5454
$ws.Range("L91").Value = 1514.5454
$ws.Range("N91").Value = -4322.5454
$ws.Range("H97").Value = 5537.8213
$ws.Range("I97").Value = 5645.8096
$ws.Range("K97").Value = 5645.8096
$ws.Range("M97").Value = -5149.8096
$ws.Range("H122").Value = 3640.7778
$ws.Range("I122").Value = 3345.875
$ws.Range("K122").Value = 10037.625
$ws.Range("M122").Value = -7587.625
$ws.Range("H132").Value = 40612.04
$ws.Range("I132").Value = 40612.04
$ws.Range("K132").Value = 121836.12
$ws.Range("M132").Value = -119306.12
$ws.Range("H136").Value = 1955
$ws.Range("I136").Value = 1786.3334
$ws.Range("K136").Value = 5359.0002
$ws.Range("M136").Value = -2809.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 31200
$ws.Range("J6").Value = 31200
$ws.Range("L6").Value = 31200
$ws.Range("N6").Value = -31426
$ws.Range("H22").Value = 167299.83
$ws.Range("I22").Value = 200719.8
$ws.Range("K22").Value = 200719.8
$ws.Range("M22").Value = -200546.8
$ws.Range("H82").Value = 9628.166999999999
$ws.Range("I82").Value = 9628.166999999999
$ws.Range("K82").Value = 9628.166999999999
$ws.Range("M82").Value = -9245.166999999999
$ws.Range("H85").Value = 9628.166999999999
$ws.Range("I85").Value = 9628.166999999999
$ws.Range("K85").Value = 9628.166999999999
$ws.Range("M85").Value = -8302.166999999999
$ws.Range("H94").Value = 3123.2856
$ws.Range("J94").Value = 3577.2
$ws.Range("L94").Value = 3577.2
$ws.Range("N94").Value = -4479.2
$ws.Range("H97").Value = 16831.75
$ws.Range("I97").Value = 13109
$ws.Range("K97").Value = 13109
$ws.Range("M97").Value = -12118
$ws.Range("H99").Value = 29933.879
$ws.Range("I99").Value = 86314.25
$ws.Range("J99").Value = 6604.069
$ws.Range("K99").Value = 86314.25
$ws.Range("L99").Value = 6604.069
$ws.Range("M99").Value = -84816.25
$ws.Range("N99").Value = -9600.069
$ws.Range("H105").Value = 5508.317
$ws.Range("I105").Value = 5482.143
$ws.Range("K105").Value = 5482.143
$ws.Range("M105").Value = -3735.143
$ws.Range("H107").Value = 1798.9
$ws.Range("I107").Value = 998.625
$ws.Range("K107").Value = 998.625
$ws.Range("M107").Value = 921.375
$ws.Range("H134").Value = 3468.1428
$ws.Range("I134").Value = 3113.6667
$ws.Range("J134").Value = 5595
$ws.Range("K134").Value = 9341.000100000001
$ws.Range("L134").Value = 16785
$ws.Range("M134").Value = -6806.000100000001
$ws.Range("N134").Value = -21855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3495
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3495
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3495
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4085
$ws.Range("H34").Value = 3495
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3495
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3495
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3899
$ws.Range("H58").Value = 69528
$ws.Range("I58").Value = 93382
$ws.Range("K58").Value = 93382
$ws.Range("M58").Value = -93179
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H86").Value = 70584.164
$ws.Range("I86").Value = 103876.25
$ws.Range("K86").Value = 103876.25
$ws.Range("M86").Value = -102753.25
$ws.Range("H89").Value = 70584.164
$ws.Range("I89").Value = 103876.25
$ws.Range("K89").Value = 519381.25
$ws.Range("M89").Value = -513765.25
$ws.Range("H105").Value = 3698.25
$ws.Range("I105").Value = 3931
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 3931
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -2184
$ws.Range("N105").Value = -6494
$ws.Range("H122").Value = 2490.8572
$ws.Range("I122").Value = 2490.8572
$ws.Range("K122").Value = 7472.571599999999
$ws.Range("M122").Value = -5022.571599999999
$ws.Range("H134").Value = 50630.668
$ws.Range("I134").Value = 73796.78999999999
$ws.Range("J134").Value = 4298.4287
$ws.Range("K134").Value = 221390.37
$ws.Range("L134").Value = 12895.2861
$ws.Range("M134").Value = -218855.37
$ws.Range("N134").Value = -17965.2861
$ws.Range("H136").Value = 69528
$ws.Range("I136").Value = 93382
$ws.Range("K136").Value = 280146
$ws.Range("M136").Value = -277596

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1508313.2
$ws.Range("I4").Value = 1344709.9
$ws.Range("K4").Value = 4034129.7
$ws.Range("M4").Value = -4034017.7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H26").Value = 186.8
$ws.Range("I26").Value = 184.75
$ws.Range("J26").Value = 195
$ws.Range("K26").Value = 554.25
$ws.Range("L26").Value = 585
$ws.Range("M26").Value = -266.25
$ws.Range("N26").Value = -1161
$ws.Range("H98").Value = 3288.6667
$ws.Range("I98").Value = 1300
$ws.Range("J98").Value = 3686.4
$ws.Range("K98").Value = 3900
$ws.Range("L98").Value = 11059.2
$ws.Range("M98").Value = -2402
$ws.Range("N98").Value = -14055.2
$ws.Range("H113").Value = 932.1818
$ws.Range("J113").Value = 902.125
$ws.Range("L113").Value = 2706.375
$ws.Range("N113").Value = -7046.375
$ws.Range("H119").Value = 12014
$ws.Range("I119").Value = 4029
$ws.Range("J119").Value = 19999
$ws.Range("K119").Value = 12087
$ws.Range("L119").Value = 59997
$ws.Range("M119").Value = -7249
$ws.Range("N119").Value = -69673
$ws.Range("H121").Value = 1288.5
$ws.Range("I121").Value = 1492
$ws.Range("J121").Value = 1026.8572
$ws.Range("K121").Value = 4476
$ws.Range("L121").Value = 3080.5716
$ws.Range("M121").Value = -3166
$ws.Range("N121").Value = -5700.571599999999
$ws.Range("H131").Value = 3716018
$ws.Range("I131").Value = 7195.4116
$ws.Range("J131").Value = 10021016
$ws.Range("K131").Value = 21586.2348
$ws.Range("L131").Value = 30063048
$ws.Range("M131").Value = -16546.2348
$ws.Range("N131").Value = -30073128
$ws.Range("H137").Value = 4148.091
$ws.Range("I137").Value = 1565
$ws.Range("J137").Value = 4722.1113
$ws.Range("K137").Value = 4695
$ws.Range("L137").Value = 14166.3339
$ws.Range("M137").Value = 405
$ws.Range("N137").Value = -24366.3339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8364.799999999999
$ws.Range("I70").Value = 7664.5293
$ws.Range("K70").Value = 7664.5293
$ws.Range("M70").Value = -7394.5293
$ws.Range("H73").Value = 8364.799999999999
$ws.Range("I73").Value = 7664.5293
$ws.Range("K73").Value = 7664.5293
$ws.Range("M73").Value = -6728.5293
$ws.Range("H80").Value = 2726.0454
$ws.Range("I80").Value = 2025.4546
$ws.Range("J80").Value = 3426.6365
$ws.Range("K80").Value = 2025.4546
$ws.Range("L80").Value = 3426.6365
$ws.Range("M80").Value = -1027.4546
$ws.Range("N80").Value = -5422.636500000001
$ws.Range("H83").Value = 2726.0454
$ws.Range("I83").Value = 2025.4546
$ws.Range("J83").Value = 3426.6365
$ws.Range("K83").Value = 10127.273
$ws.Range("L83").Value = 17133.1825
$ws.Range("M83").Value = -5135.273000000001
$ws.Range("N83").Value = -27117.1825
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H97").Value = 1098.5
$ws.Range("I97").Value = 829.375
$ws.Range("J97").Value = 2175
$ws.Range("K97").Value = 829.375
$ws.Range("L97").Value = 2175
$ws.Range("M97").Value = -333.375
$ws.Range("N97").Value = -3167
$ws.Range("H102").Value = 2060.577
$ws.Range("I102").Value = 1546.7826
$ws.Range("K102").Value = 1546.7826
$ws.Range("M102").Value = 75.2174
$ws.Range("H117").Value = 64000
$ws.Range("J117").Value = 64000
$ws.Range("L117").Value = 64000
$ws.Range("N117").Value = -70884
$ws.Range("H122").Value = 2989.4285
$ws.Range("I122").Value = 2529.6
$ws.Range("J122").Value = 3244.889
$ws.Range("K122").Value = 7588.799999999999
$ws.Range("L122").Value = 9734.667000000001
$ws.Range("M122").Value = -5138.799999999999
$ws.Range("N122").Value = -14634.667
$ws.Range("H132").Value = 93948.91
$ws.Range("I132").Value = 114070.89
$ws.Range("K132").Value = 342212.67
$ws.Range("M132").Value = -339682.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 67223.3
$ws.Range("J22").Value = 4932.3335
$ws.Range("L22").Value = 4932.3335
$ws.Range("N22").Value = -5522.3335
$ws.Range("H27").Value = 67223.3
$ws.Range("J27").Value = 4932.3335
$ws.Range("L27").Value = 4932.3335
$ws.Range("N27").Value = -5146.3335
$ws.Range("H46").Value = 16846.572
$ws.Range("I46").Value = 21818.6
$ws.Range("J46").Value = 4416.5
$ws.Range("K46").Value = 21818.6
$ws.Range("L46").Value = 4416.5
$ws.Range("M46").Value = -21630.6
$ws.Range("N46").Value = -4792.5
$ws.Range("H55").Value = 834.375
$ws.Range("I55").Value = 739.2857
$ws.Range("J55").Value = 1500
$ws.Range("K55").Value = 739.2857
$ws.Range("L55").Value = 1500
$ws.Range("M55").Value = -566.2857
$ws.Range("N55").Value = -1846
$ws.Range("H74").Value = 62249.25
$ws.Range("I74").Value = 54000
$ws.Range("K74").Value = 54000
$ws.Range("M74").Value = -53002
$ws.Range("H77").Value = 62249.25
$ws.Range("I77").Value = 54000
$ws.Range("K77").Value = 162000
$ws.Range("M77").Value = -157008
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H122").Value = 3867.2258
$ws.Range("I122").Value = 3232.6667
$ws.Range("J122").Value = 4745.846
$ws.Range("K122").Value = 9698.000100000001
$ws.Range("L122").Value = 14237.538
$ws.Range("M122").Value = -7248.000100000001
$ws.Range("N122").Value = -19137.538
$ws.Range("H132").Value = 47450.184
$ws.Range("I132").Value = 54724.39
$ws.Range("K132").Value = 164173.17
$ws.Range("M132").Value = -161643.17

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 53782
$ws.Range("I45").Value = 9599
$ws.Range("J45").Value = 64827.75
$ws.Range("K45").Value = 9599
$ws.Range("L45").Value = 64827.75
$ws.Range("M45").Value = -9108
$ws.Range("N45").Value = -65809.75
$ws.Range("H81").Value = 1584.5454
$ws.Range("I81").Value = 1343
$ws.Range("K81").Value = 2686
$ws.Range("M81").Value = -1625
$ws.Range("H84").Value = 1584.5454
$ws.Range("I84").Value = 1343
$ws.Range("K84").Value = 13430
$ws.Range("M84").Value = -8126
$ws.Range("H122").Value = 8445.049999999999
$ws.Range("I122").Value = 9241.117
$ws.Range("K122").Value = 27723.351
$ws.Range("M122").Value = -25273.351
$ws.Range("H132").Value = 240365.25
$ws.Range("I132").Value = 261307.55
$ws.Range("K132").Value = 783922.6499999999
$ws.Range("M132").Value = -781392.6499999999
